$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("case")
$ws.Rows.Item(60).Delete()
